$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.84"
$ws.Range("E2").Value = "'-1.25%"
$ws.Range("D3").Value = "'44.38"
$ws.Range("E3").Value = "'2.68%"
$ws.Range("D4").Value = "'5.577"
$ws.Range("E4").Value = "'-1.73%"
$ws.Range("D5").Value = "'0.08060"
$ws.Range("E5").Value = "'-1.99%"
$ws.Range("D6").Value = "'4.303"
$ws.Range("E6").Value = "'-4.99%"
$ws.Range("D7").Value = "'1.899"
$ws.Range("E7").Value = "'-1.98%"
$ws.Range("D8").Value = "'2.610"
$ws.Range("E8").Value = "'-8.57%"
$ws.Range("D9").Value = "'0.9468"
$ws.Range("E9").Value = "'0.11%"
$ws.Range("D10").Value = "'0.1166"
$ws.Range("E10").Value = "'-3.94%"
$ws.Range("D11").Value = "'0.1833"
$ws.Range("E11").Value = "'-5.86%"
$ws.Range("D12").Value = "'0.09653"
$ws.Range("E12").Value = "'-1.11%"
$ws.Range("D13").Value = "'0.04370"
$ws.Range("E13").Value = "'-4.25%"
$ws.Range("D14").Value = "'0.1068"
$ws.Range("E14").Value = "'0.18%"
$ws.Range("D15").Value = "'0.001284"
$ws.Range("E15").Value = "'0.48%"
$ws.Range("D16").Value = "'0.005899"
$ws.Range("E16").Value = "'-2.28%"
$ws.Range("D17").Value = "'3.628"
$ws.Range("E17").Value = "'4.26%"
$ws.Range("D18").Value = "'0.3500"
$ws.Range("E18").Value = "'-1.05%"
$ws.Range("D19").Value = "'9.604"
$ws.Range("E19").Value = "'9.12%"
$ws.Range("D20").Value = "'0.1388"
$ws.Range("E20").Value = "'2.76%"
$ws.Range("D21").Value = "'0.2651"
$ws.Range("E21").Value = "'-2.60%"
$ws.Range("D22").Value = "'0.04223"
$ws.Range("E22").Value = "'-4.29%"
$ws.Range("D23").Value = "'0.001247"
$ws.Range("E23").Value = "'0.16%"
$ws.Range("D24").Value = "'0.004510"
$ws.Range("E24").Value = "'4.12%"
$ws.Range("D25").Value = "'0.0001261"
$ws.Range("E25").Value = "'1.92%"
$ws.Range("D26").Value = "'0.0003992"
$ws.Range("E26").Value = "'-0.55%"
$ws.Range("D38").Value = "'0.02645"
$ws.Range("E38").Value = "'-6.00%"
$ws.Range("D39").Value = "'0.05503"
$ws.Range("E39").Value = "'-4.17%"
$ws.Range("D40").Value = "'0.007590"
$ws.Range("E40").Value = "'-4.27%"
$ws.Range("D41").Value = "'0.1400"
$ws.Range("E41").Value = "'-0.68%"
$ws.Range("D42").Value = "'0.007074"
$ws.Range("E42").Value = "'-24.94%"
$ws.Range("E43").Value = "'-6.37%"
$ws.Range("D44").Value = "'0.008367"
$ws.Range("E44").Value = "'-15.63%"
$ws.Range("D45").Value = "'0.00006922"
$ws.Range("E45").Value = "'-9.37%"
$ws.Range("E46").Value = "'-0.52%"
$ws.Range("D47").Value = "'0.002272"
$ws.Range("E47").Value = "'-0.55%"
$ws.Range("D48").Value = "'0.003838"
$ws.Range("E48").Value = "'20.29%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'-0.52%"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'-0.52%"
